$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows entirely (rather than just clearing their contents)
$null = $ws.Rows("2:3").Delete()

# Update A1 with the new contact name
$ws.Range("A1").Value = "Julio César Pérez"

# Make sure A1 is the selected/active cell, matching the saved view state
$null = $ws.Range("A1").Select()
